# Update data.xlsx from the QR tool export
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("locations")

# Remove the old "Test 3WEL" demo row (row 4). This shifts the
# "Tran Van A" template row (old row 5) up into row 4.
$ws.Rows(4).Delete()

# --- Row 2: Cá nhân Huỳnh Thị Thanh Thúy -------------------------------
$ws.Cells.Item(2, 1).Value = "8p1jaw70gl3"
$ws.Cells.Item(2, 2).Value = "m396pkz6"
$ws.Cells.Item(2, 3).Value = "Cá nhân Huỳnh Thị Thanh Thúy"
$ws.Cells.Item(2, 6).Value = "2025-08-14T06:34:09.443Z"
$ws.Cells.Item(2, 7).Value = "Nhà mặt tiền giá trị cao"
$ws.Cells.Item(2, 8).Value = "09134563456"
$ws.Cells.Item(2, 10).Value = "123456789102"
$ws.Cells.Item(2, 11).Value = "1890000000"
$ws.Cells.Item(2, 12).Value = "Huỳnh Thị Thanh Thúy"
$ws.Cells.Item(2, 13).Value = "40304b41a6114ed9"
$ws.Cells.Item(2, 14).Value = "bde53dff6af6c2afebd7bcbd1b423ae3c27e01a566782dd87b0ae00ad65d259b"

# --- Row 3: Hộ kinh doanh Test TFS3 -------------------------------------
$ws.Cells.Item(3, 1).Value = "q7aklukam98"
$ws.Cells.Item(3, 2).Value = "TESTXO03"
$ws.Cells.Item(3, 3).Value = "Hộ kinh doanh Test TFS3"
$ws.Cells.Item(3, 6).Value = "2025-08-14T06:21:07.541Z"
$ws.Cells.Item(3, 8).Value = "0908699201"
$ws.Cells.Item(3, 10).Value = "0339828908379"
$ws.Cells.Item(3, 11).Value = "KH5502"

# --- Row 4 (was row 5): Hộ kinh doanh Trần Văn A ------------------------
$ws.Cells.Item(4, 1).Value = "lt7ouixns1"
$ws.Cells.Item(4, 2).Value = "DEMOJGGQ"
$ws.Cells.Item(4, 6).Value = "2025-08-14T06:18:53.958Z"
$ws.Cells.Item(4, 7).Value = "Dòng mẫu để thử"
$ws.Cells.Item(4, 9).Value = "CN Cần Thơ II"
